$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Insert two new columns before AK (former AK:AP data shifts right to AM:AR)
$ws.Range("AK1:AL1").EntireColumn.Insert()

# New header labels for the inserted columns (row 1 keeps the header styling
# that was already copied onto AK1/AL1 by the column insert)
$ws.Range("AK1").Value = "measure11_highPerf"
$ws.Range("AL1").Value = "measure11_lowPerf"

# Give the new columns the same width used by nearby wide columns (e.g. AI)
$ws.Range("AK1:AL1").ColumnWidth = 24

# Add a new row 12 below the table; C12 reuses the date-format style from the
# C column (e.g. C2) but stays empty, matching a "reset row" placeholder.
$ws.Range("C2").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C12").ClearContents()

# Leave the selection where the author ended up after the edit
$ws.Range("AK14").Select() | Out-Null
